$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 155.375
$ws.Range("I42").Value = 64
$ws.Range("J42").Value = 185.83333
$ws.Range("K42").Value = 192
$ws.Range("L42").Value = 557.49999
$ws.Range("M42").Value = 38
$ws.Range("N42").Value = -1017.49999
$ws.Range("H53").Value = 612.7308
$ws.Range("I53").Value = 237.72223
$ws.Range("J53").Value = 1456.5
$ws.Range("K53").Value = 237.72223
$ws.Range("L53").Value = 1456.5
$ws.Range("M53").Value = 399.27777
$ws.Range("N53").Value = -2730.5
$ws.Range("H82").Value = 16667733
$ws.Range("I82").Value = 16667733
$ws.Range("K82").Value = 50003199
$ws.Range("M82").Value = -50002793
$ws.Range("H85").Value = 16667733
$ws.Range("I85").Value = 16667733
$ws.Range("K85").Value = 50003199
$ws.Range("M85").Value = -50001795
$ws.Range("H98").Value = 4061.8333
$ws.Range("I98").Value = 3638.3928
$ws.Range("K98").Value = 3638.3928
$ws.Range("M98").Value = -2140.3928
$ws.Range("H100").Value = 1585.8889
$ws.Range("I100").Value = 1438.4375
$ws.Range("J100").Value = 1800.3636
$ws.Range("K100").Value = 1438.4375
$ws.Range("L100").Value = 1800.3636
$ws.Range("M100").Value = -897.4375
$ws.Range("N100").Value = -2882.3636
$ws.Range("H122").Value = 4061.8333
$ws.Range("I122").Value = 3638.3928
$ws.Range("K122").Value = 10915.1784
$ws.Range("M122").Value = -8465.178400000001
$ws.Range("H127").Value = 1297.4783
$ws.Range("I127").Value = 390.33334
$ws.Range("J127").Value = 1617.6471
$ws.Range("K127").Value = 1171.00002
$ws.Range("L127").Value = 4852.9413
$ws.Range("M127").Value = 3788.99998
$ws.Range("N127").Value = -14772.9413

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4204365.5
$ws.Range("I2").Value = 2433.3333
$ws.Range("J2").Value = 7355815
$ws.Range("K2").Value = 2433.3333
$ws.Range("L2").Value = 7355815
$ws.Range("M2").Value = -2320.3333
$ws.Range("N2").Value = -7356041
$ws.Range("H102").Value = 2199.9092
$ws.Range("I102").Value = 2169.9
$ws.Range("K102").Value = 2169.9
$ws.Range("M102").Value = -547.9000000000001
$ws.Range("H116").Value = 4204365.5
$ws.Range("I116").Value = 2433.3333
$ws.Range("J116").Value = 7355815
$ws.Range("K116").Value = 2433.3333
$ws.Range("L116").Value = 7355815
$ws.Range("M116").Value = -139.3332999999998
$ws.Range("N116").Value = -7360403
$ws.Range("H122").Value = 1935.7916
$ws.Range("I122").Value = 1684.9375
$ws.Range("J122").Value = 2437.5
$ws.Range("K122").Value = 5054.8125
$ws.Range("L122").Value = 7312.5
$ws.Range("M122").Value = -2604.8125
$ws.Range("N122").Value = -12212.5
$ws.Range("H139").Value = 54828.57
$ws.Range("J139").Value = 54828.57
$ws.Range("L139").Value = 54828.57
$ws.Range("N139").Value = -65108.57

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4204365.5
$ws.Range("I3").Value = 2433.3333
$ws.Range("J3").Value = 7355815
$ws.Range("K3").Value = 2433.3333
$ws.Range("L3").Value = 7355815
$ws.Range("M3").Value = -2319.3333
$ws.Range("N3").Value = -7356043
$ws.Range("H86").Value = 1483.8334
$ws.Range("I86").Value = 1474.3158
$ws.Range("J86").Value = 1520
$ws.Range("K86").Value = 1474.3158
$ws.Range("L86").Value = 1520
$ws.Range("M86").Value = -351.3158000000001
$ws.Range("N86").Value = -3766
$ws.Range("H89").Value = 1483.8334
$ws.Range("I89").Value = 1474.3158
$ws.Range("J89").Value = 1520
$ws.Range("K89").Value = 7371.579000000001
$ws.Range("L89").Value = 7600
$ws.Range("M89").Value = -1755.579000000001
$ws.Range("N89").Value = -18832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 100
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 100
$ws.Range("L25").Value = 100
$ws.Range("M25").Value = 74
$ws.Range("N25").Value = -448
$ws.Range("H130").Value = 51708.57
$ws.Range("J130").Value = 51708.57
$ws.Range("L130").Value = 51708.57
$ws.Range("N130").Value = -61748.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2657.4138
$ws.Range("I131").Value = 20086.666
$ws.Range("J131").Value = 1706.7273
$ws.Range("K131").Value = 60259.99800000001
$ws.Range("L131").Value = 5120.1819
$ws.Range("M131").Value = -55219.99800000001
$ws.Range("N131").Value = -15200.1819
$ws.Range("H132").Value = 62501584
$ws.Range("I132").Value = 111112170
$ws.Range("J132").Value = 2268.2856
$ws.Range("K132").Value = 1000009530
$ws.Range("L132").Value = 20414.5704
$ws.Range("M132").Value = -1000007000
$ws.Range("N132").Value = -25474.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 59088084
$ws.Range("I132").Value = 111112616
$ws.Range("J132").Value = 2334049.2
$ws.Range("K132").Value = 333337848
$ws.Range("L132").Value = 7002147.600000001
$ws.Range("M132").Value = -333335318
$ws.Range("N132").Value = -7007207.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1147.3684
$ws.Range("I46").Value = 1033.3334
$ws.Range("J46").Value = 1575
$ws.Range("K46").Value = 1033.3334
$ws.Range("L46").Value = 1575
$ws.Range("M46").Value = -845.3334
$ws.Range("N46").Value = -1951
$ws.Range("H132").Value = 8091.409
$ws.Range("I132").Value = 16533.111
$ws.Range("J132").Value = 2247.1538
$ws.Range("K132").Value = 49599.333
$ws.Range("L132").Value = 6741.4614
$ws.Range("M132").Value = -47069.333
$ws.Range("N132").Value = -11801.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4900.2856
$ws.Range("I62").Value = 4950.5
$ws.Range("J62").Value = 4833.3335
$ws.Range("K62").Value = 4950.5
$ws.Range("L62").Value = 4833.3335
$ws.Range("M62").Value = -4326.5
$ws.Range("N62").Value = -6081.3335
$ws.Range("H65").Value = 4900.2856
$ws.Range("I65").Value = 4950.5
$ws.Range("J65").Value = 4833.3335
$ws.Range("K65").Value = 24752.5
$ws.Range("L65").Value = 24166.6675
$ws.Range("M65").Value = -21632.5
$ws.Range("N65").Value = -30406.6675
$ws.Range("H132").Value = 18822990
$ws.Range("I132").Value = 15171334
$ws.Range("J132").Value = 40002600
$ws.Range("K132").Value = 45514002
$ws.Range("L132").Value = 120007800
$ws.Range("M132").Value = -45511472
$ws.Range("N132").Value = -120012860
